# Applies the NATMI TPM data refresh for the Fgf15-Fgfr3 LR-pair sheet:
#  - Rows 2-7 get updated TPM-derived statistics (and column B/C/D cluster
#    labels are refreshed to reflect the new shared-string ordering).
#  - Rows 8-13 are newly added for the "Resolving-Mac" target cluster pairing
#    (MuSCs/Resolving-Mac senders x ECs/FAPs/MuSCs targets).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ row = 2; A = "ECs"; B = "Fgf15"; C = "Fgfr3"; D = "ECs"; E = 1; F = 0.3333333333333333; G = 0.06665; H = 0.19995; I = 0.3040875017679506; J = 0.3040875017679506; K = 3; L = 1; M = 6.045145666666667; N = 18.135437; O = 0.8160840232643366; P = 0.8160840232643367; Q = 0.4029089586833334; R = 3.626180628150001; S = 0.2481609518671902; T = 0.2481609518671902 },
    @{ row = 3; A = "ECs"; B = "Fgf15"; C = "Fgfr3"; D = "FAPs"; E = 1; F = 0.3333333333333333; G = 0.06665; H = 0.19995; I = 0.3040875017679506; J = 0.3040875017679506; K = 3; L = 1; M = 0.6824433333333334; N = 2.04733; O = 0.09212864864242169; P = 0.09212864864242169; Q = 0.04548484816666667; R = 0.4093636335; S = 0.02801517060693131; T = 0.0280151706069313 },
    @{ row = 4; A = "ECs"; B = "Fgf15"; C = "Fgfr3"; D = "MuSCs"; E = 1; F = 0.3333333333333333; G = 0.06665; H = 0.19995; I = 0.3040875017679506; J = 0.3040875017679506; K = 3; L = 1; M = 0.6799149999999999; N = 2.039745; O = 0.09178732809324164; P = 0.09178732809324165; Q = 0.04531633475; R = 0.40784701275; S = 0.02791137929382908; T = 0.02791137929382908 },
    @{ row = 5; A = "FAPs"; B = "Fgf15"; C = "Fgfr3"; D = "ECs"; E = 1; F = 0.3333333333333333; G = 0.06149566666666667; H = 0.184487; I = 0.2805710974676865; J = 0.2805710974676864; K = 3; L = 1; M = 6.045145666666667; N = 18.135437; O = 0.8160840232643366; P = 0.8160840232643367; Q = 0.3717502628687778; R = 3.345752365819001; S = 0.2289695900331199; T = 0.2289695900331199 },
    @{ row = 6; A = "FAPs"; B = "Fgf15"; C = "Fgfr3"; D = "FAPs"; E = 1; F = 0.3333333333333333; G = 0.06149566666666667; H = 0.184487; I = 0.2805710974676865; J = 0.2805710974676864; K = 3; L = 1; M = 0.6824433333333334; N = 2.04733; O = 0.09212864864242169; P = 0.09212864864242169; Q = 0.04196730774555556; R = 0.37770576971; S = 0.02584863605781914; T = 0.02584863605781914 },
    @{ row = 7; A = "FAPs"; B = "Fgf15"; C = "Fgfr3"; D = "MuSCs"; E = 1; F = 0.3333333333333333; G = 0.06149566666666667; H = 0.184487; I = 0.2805710974676865; J = 0.2805710974676864; K = 3; L = 1; M = 0.6799149999999999; N = 2.039745; O = 0.09178732809324164; P = 0.09178732809324165; Q = 0.04181182620166667; R = 0.376306435815; S = 0.02575287137674742; T = 0.02575287137674742 },
    @{ row = 8; A = "MuSCs"; B = "Fgf15"; C = "Fgfr3"; D = "ECs"; E = 2; F = 0.6666666666666666; G = 0.057558; H = 0.172674; I = 0.2626056778208508; J = 0.2626056778208507; K = 3; L = 1; M = 6.045145666666667; N = 18.135437; O = 0.8160840232643366; P = 0.8160840232643367; Q = 0.347946494282; R = 3.131518448538; S = 0.2143082980880981; T = 0.214308298088098 },
    @{ row = 9; A = "MuSCs"; B = "Fgf15"; C = "Fgfr3"; D = "FAPs"; E = 2; F = 0.6666666666666666; G = 0.057558; H = 0.172674; I = 0.2626056778208508; J = 0.2626056778208507; K = 3; L = 1; M = 0.6824433333333334; N = 2.04733; O = 0.09212864864242169; P = 0.09212864864242169; Q = 0.03928007338; R = 0.35352066042; S = 0.02419350622346215; T = 0.02419350622346214 },
    @{ row = 10; A = "MuSCs"; B = "Fgf15"; C = "Fgfr3"; D = "MuSCs"; E = 2; F = 0.6666666666666666; G = 0.057558; H = 0.172674; I = 0.2626056778208508; J = 0.2626056778208507; K = 3; L = 1; M = 0.6799149999999999; N = 2.039745; O = 0.09178732809324164; P = 0.09178732809324165; Q = 0.03913454757; R = 0.35221092813; S = 0.02410387350929054; T = 0.02410387350929053 },
    @{ row = 11; A = "Resolving-Mac"; B = "Fgf15"; C = "Fgfr3"; D = "ECs"; E = 2; F = 0.6666666666666666; G = 0.03347666666666666; H = 0.10043; I = 0.1527357229435123; J = 0.1527357229435123; K = 3; L = 1; M = 6.045145666666667; N = 18.135437; O = 0.8160840232643366; P = 0.8160840232643367; Q = 0.2023713264344444; R = 1.82134193791; S = 0.1246451832759285; T = 0.1246451832759286 },
    @{ row = 12; A = "Resolving-Mac"; B = "Fgf15"; C = "Fgfr3"; D = "FAPs"; E = 2; F = 0.6666666666666666; G = 0.03347666666666666; H = 0.10043; I = 0.1527357229435123; J = 0.1527357229435123; K = 3; L = 1; M = 0.6824433333333334; N = 2.04733; O = 0.09212864864242169; P = 0.09212864864242169; Q = 0.02284592798888889; R = 0.2056133519; S = 0.01407133575420911; T = 0.01407133575420911 },
    @{ row = 13; A = "Resolving-Mac"; B = "Fgf15"; C = "Fgfr3"; D = "MuSCs"; E = 2; F = 0.6666666666666666; G = 0.03347666666666666; H = 0.10043; I = 0.1527357229435123; J = 0.1527357229435123; K = 3; L = 1; M = 0.6799149999999999; N = 2.039745; O = 0.09178732809324164; P = 0.09178732809324165; Q = 0.02276128781666666; R = 0.20485159035; S = 0.01401920391337461; T = 0.01401920391337462 }
)

foreach ($r in $rowData) {
    if ($r.ContainsKey("A")) { $ws.Range("A" + $r.row).Value = $r.A }
    if ($r.ContainsKey("B")) { $ws.Range("B" + $r.row).Value = $r.B }
    if ($r.ContainsKey("C")) { $ws.Range("C" + $r.row).Value = $r.C }
    if ($r.ContainsKey("D")) { $ws.Range("D" + $r.row).Value = $r.D }
    if ($r.ContainsKey("E")) { $ws.Range("E" + $r.row).Value = $r.E }
    if ($r.ContainsKey("F")) { $ws.Range("F" + $r.row).Value = $r.F }
    if ($r.ContainsKey("G")) { $ws.Range("G" + $r.row).Value = $r.G }
    if ($r.ContainsKey("H")) { $ws.Range("H" + $r.row).Value = $r.H }
    if ($r.ContainsKey("I")) { $ws.Range("I" + $r.row).Value = $r.I }
    if ($r.ContainsKey("J")) { $ws.Range("J" + $r.row).Value = $r.J }
    if ($r.ContainsKey("K")) { $ws.Range("K" + $r.row).Value = $r.K }
    if ($r.ContainsKey("L")) { $ws.Range("L" + $r.row).Value = $r.L }
    if ($r.ContainsKey("M")) { $ws.Range("M" + $r.row).Value = $r.M }
    if ($r.ContainsKey("N")) { $ws.Range("N" + $r.row).Value = $r.N }
    if ($r.ContainsKey("O")) { $ws.Range("O" + $r.row).Value = $r.O }
    if ($r.ContainsKey("P")) { $ws.Range("P" + $r.row).Value = $r.P }
    if ($r.ContainsKey("Q")) { $ws.Range("Q" + $r.row).Value = $r.Q }
    if ($r.ContainsKey("R")) { $ws.Range("R" + $r.row).Value = $r.R }
    if ($r.ContainsKey("S")) { $ws.Range("S" + $r.row).Value = $r.S }
    if ($r.ContainsKey("T")) { $ws.Range("T" + $r.row).Value = $r.T }
}
